$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 397.42856
$ws.Range("I6").Value = 440
$ws.Range("J6").Value = 291
$ws.Range("K6").Value = 1320
$ws.Range("L6").Value = 873
$ws.Range("M6").Value = -1208
$ws.Range("N6").Value = -1097

$ws.Range("H111").Value = 4163.7144
$ws.Range("I111").Value = 4664
$ws.Range("J111").Value = 3496.6667
$ws.Range("K111").Value = 13992
$ws.Range("L111").Value = 10490.0001
$ws.Range("M111").Value = -10925
$ws.Range("N111").Value = -16624.0001

$ws.Range("H125").Value = 5179.121
$ws.Range("I125").Value = 7801.9443
$ws.Range("J125").Value = 2031.7333
$ws.Range("K125").Value = 70217.4987
$ws.Range("L125").Value = 18285.5997
$ws.Range("M125").Value = -67757.4987
$ws.Range("N125").Value = -23205.5997

$ws.Range("H135").Value = 885.84375
$ws.Range("I135").Value = 417.8
$ws.Range("J135").Value = 1665.9166
$ws.Range("K135").Value = 3760.2
$ws.Range("L135").Value = 14993.2494
$ws.Range("M135").Value = -1225.2
$ws.Range("N135").Value = -20063.2494

$ws.Range("H138").Value = 2904.0857
$ws.Range("I138").Value = 2259.65
$ws.Range("J138").Value = 3763.3333
$ws.Range("K138").Value = 6778.950000000001
$ws.Range("L138").Value = 11289.9999
$ws.Range("M138").Value = -1638.950000000001
$ws.Range("N138").Value = -21569.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2713.182
$ws.Range("I32").Value = 2539.157
$ws.Range("J32").Value = 4932
$ws.Range("K32").Value = 2539.157
$ws.Range("L32").Value = 4932
$ws.Range("M32").Value = -2252.157
$ws.Range("N32").Value = -5506

$ws.Range("H61").Value = 23811248
$ws.Range("I61").Value = 30304558
$ws.Range("J61").Value = 2441.111
$ws.Range("K61").Value = 30304558
$ws.Range("L61").Value = 2441.111
$ws.Range("M61").Value = -30304346
$ws.Range("N61").Value = -2865.111

$ws.Range("H132").Value = 5342.569
$ws.Range("I132").Value = 4089.244
$ws.Range("J132").Value = 8365.294
$ws.Range("K132").Value = 12267.732
$ws.Range("L132").Value = 25095.882
$ws.Range("M132").Value = -9737.732
$ws.Range("N132").Value = -30155.882

$ws.Range("H136").Value = 23811248
$ws.Range("I136").Value = 30304558
$ws.Range("J136").Value = 2441.111
$ws.Range("K136").Value = 90913674
$ws.Range("L136").Value = 7323.333
$ws.Range("M136").Value = -90911124
$ws.Range("N136").Value = -12423.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1853.3334
$ws.Range("I134").Value = 1630.6
$ws.Range("J134").Value = 2489.7144
$ws.Range("K134").Value = 4891.799999999999
$ws.Range("L134").Value = 7469.1432
$ws.Range("M134").Value = -2356.799999999999
$ws.Range("N134").Value = -12539.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4860.6294
$ws.Range("I31").Value = 1799.8462
$ws.Range("K31").Value = 1799.8462
$ws.Range("M31").Value = -1504.8462

$ws.Range("H34").Value = 4860.6294
$ws.Range("I34").Value = 1799.8462
$ws.Range("K34").Value = 1799.8462
$ws.Range("M34").Value = -1597.8462

$ws.Range("H58").Value = 2099.5
$ws.Range("I58").Value = 2418.8572
$ws.Range("J58").Value = 1927.5385
$ws.Range("K58").Value = 2418.8572
$ws.Range("L58").Value = 1927.5385
$ws.Range("M58").Value = -2215.8572
$ws.Range("N58").Value = -2333.5385

$ws.Range("H132").Value = 55564760
$ws.Range("I132").Value = 100013570
$ws.Range("J132").Value = 3749.75
$ws.Range("K132").Value = 300040710
$ws.Range("L132").Value = 11249.25
$ws.Range("M132").Value = -300038180
$ws.Range("N132").Value = -16309.25

$ws.Range("H134").Value = 8221.6
$ws.Range("I134").Value = 16240
$ws.Range("J134").Value = 2876
$ws.Range("K134").Value = 48720
$ws.Range("L134").Value = 8628
$ws.Range("M134").Value = -46185
$ws.Range("N134").Value = -13698

$ws.Range("H136").Value = 2099.5
$ws.Range("I136").Value = 2418.8572
$ws.Range("J136").Value = 1927.5385
$ws.Range("K136").Value = 7256.571599999999
$ws.Range("L136").Value = 5782.6155
$ws.Range("M136").Value = -4706.571599999999
$ws.Range("N136").Value = -10882.6155

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 990.2929
$ws.Range("I68").Value = 745.6875
$ws.Range("J68").Value = 1107.1194
$ws.Range("K68").Value = 2237.0625
$ws.Range("L68").Value = 3321.3582
$ws.Range("M68").Value = -1426.0625
$ws.Range("N68").Value = -4943.358200000001

$ws.Range("H71").Value = 990.2929
$ws.Range("I71").Value = 745.6875
$ws.Range("J71").Value = 1107.1194
$ws.Range("K71").Value = 6711.1875
$ws.Range("L71").Value = 9964.0746
$ws.Range("M71").Value = -2655.1875
$ws.Range("N71").Value = -18076.0746

$ws.Range("H97").Value = 358.5
$ws.Range("I97").Value = 150
$ws.Range("J97").Value = 428
$ws.Range("K97").Value = 450
$ws.Range("L97").Value = 1284
$ws.Range("M97").Value = 46
$ws.Range("N97").Value = -2276

$ws.Range("H98").Value = 579.375
$ws.Range("I98").Value = 335
$ws.Range("J98").Value = 1068.125
$ws.Range("K98").Value = 1005
$ws.Range("L98").Value = 3204.375
$ws.Range("M98").Value = 493
$ws.Range("N98").Value = -6200.375

$ws.Range("H131").Value = 391.38202
$ws.Range("I131").Value = 173.17188
$ws.Range("J131").Value = 950
$ws.Range("K131").Value = 519.51564
$ws.Range("L131").Value = 2850
$ws.Range("M131").Value = 4520.48436
$ws.Range("N131").Value = -12930

$ws.Range("H134").Value = 45092.63
$ws.Range("I134").Value = 66353.7
$ws.Range("J134").Value = 8948.799999999999
$ws.Range("K134").Value = 199061.1
$ws.Range("L134").Value = 26846.4
$ws.Range("M134").Value = -193991.1
$ws.Range("N134").Value = -36986.39999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 697.5454999999999
$ws.Range("I97").Value = 723.9474
$ws.Range("J97").Value = 530.3333
$ws.Range("K97").Value = 723.9474
$ws.Range("L97").Value = 530.3333
$ws.Range("M97").Value = -227.9474
$ws.Range("N97").Value = -1522.3333

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").Value = ""

$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").Value = ""

$ws.Range("H132").Value = 4237.8887
$ws.Range("I132").Value = 4238.6484
$ws.Range("J132").Value = 4234.375
$ws.Range("K132").Value = 12715.9452
$ws.Range("L132").Value = 12703.125
$ws.Range("M132").Value = -10185.9452
$ws.Range("N132").Value = -17763.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3191
$ws.Range("I93").Value = 3083.9443
$ws.Range("J93").Value = 3833.3333
$ws.Range("K93").Value = 3083.9443
$ws.Range("L93").Value = 3833.3333
$ws.Range("M93").Value = -1835.9443
$ws.Range("N93").Value = -6329.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 43500
$ws.Range("J133").Value = 43500
$ws.Range("L133").Value = 43500
$ws.Range("N133").Value = -53620

$ws.Range("H136").Value = 6260.696
$ws.Range("I136").Value = 8924.308000000001
$ws.Range("J136").Value = 2798
$ws.Range("K136").Value = 26772.924
$ws.Range("L136").Value = 8394
$ws.Range("M136").Value = -24222.924
$ws.Range("N136").Value = -13494
